$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 8391.77  # H33: was 6958.7334
$ws.Cells.Item(33, 9).Value = 127  # I33: was 114
$ws.Cells.Item(33, 10).Value = 9894.454  # J33: was 8669.916999999999
$ws.Cells.Item(33, 11).Value = 127  # K33: was 114
$ws.Cells.Item(33, 12).Value = 9894.454  # L33: was 8669.916999999999
$ws.Cells.Item(33, 13).Value = 102  # M33: was 115
$ws.Cells.Item(33, 14).Value = -10352.454  # N33: was -9127.916999999999

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1600.3334  # H70: was 1675.375
$ws.Cells.Item(70, 9).Value = 1000  # I70: was 0
$ws.Cells.Item(70, 10).Value = 1771.8572  # J70: was 1675.375
$ws.Cells.Item(70, 11).Value = 3000  # K70: was 0
$ws.Cells.Item(70, 12).Value = 5315.571599999999  # L70: was 5026.125
$ws.Cells.Item(70, 13).Value = -2730  # M70: was None
$ws.Cells.Item(70, 14).Value = -5855.571599999999  # N70: was -5566.125

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 1600.3334  # H73: was 1675.375
$ws.Cells.Item(73, 9).Value = 1000  # I73: was 0
$ws.Cells.Item(73, 10).Value = 1771.8572  # J73: was 1675.375
$ws.Cells.Item(73, 11).Value = 3000  # K73: was 0
$ws.Cells.Item(73, 12).Value = 5315.571599999999  # L73: was 5026.125
$ws.Cells.Item(73, 13).Value = -2064  # M73: was None
$ws.Cells.Item(73, 14).Value = -7187.571599999999  # N73: was -6898.125

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2372.1333  # H113: was 2352.5386
$ws.Cells.Item(113, 9).Value = 1925.6364  # I113: was 1968.3
$ws.Cells.Item(113, 10).Value = 3600  # J113: was 3633.3333
$ws.Cells.Item(113, 11).Value = 1925.6364  # K113: was 1968.3
$ws.Cells.Item(113, 12).Value = 3600  # L113: was 3633.3333
$ws.Cells.Item(113, 13).Value = 1328.3636  # M113: was 1285.7
$ws.Cells.Item(113, 14).Value = -10108  # N113: was -10141.3333

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2245.7058  # H116: was 2443.1428
$ws.Cells.Item(116, 9).Value = 1945.1333  # I116: was 2100.3333
$ws.Cells.Item(116, 11).Value = 1945.1333  # K116: was 2100.3333
$ws.Cells.Item(116, 13).Value = 1496.8667  # M116: was 1341.6667

# ALC row 118
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(118, 8).Value = 1025.75  # H118: was 1009.0909
$ws.Cells.Item(118, 10).Value = 2141.8  # J118: was 2375
$ws.Cells.Item(118, 12).Value = 6425.400000000001  # L118: was 7125
$ws.Cells.Item(118, 14).Value = -9739.400000000001  # N118: was -10439

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1027.9722  # H132: was 989.46155
$ws.Cells.Item(132, 9).Value = 810.7  # I132: was 783.32355
$ws.Cells.Item(132, 10).Value = 2114.3333  # J132: was 2391.2
$ws.Cells.Item(132, 11).Value = 2432.1  # K132: was 2349.97065
$ws.Cells.Item(132, 12).Value = 6342.999899999999  # L132: was 7173.599999999999
$ws.Cells.Item(132, 13).Value = 97.89999999999964  # M132: was 180.0293500000002
$ws.Cells.Item(132, 14).Value = -11402.9999  # N132: was -12233.6

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1476.921  # H2: was 1483.3948
$ws.Cells.Item(2, 9).Value = 1073.7142  # I2: was 1062.862
$ws.Cells.Item(2, 10).Value = 2605.9  # J2: was 2838.4443
$ws.Cells.Item(2, 11).Value = 1073.7142  # K2: was 1062.862
$ws.Cells.Item(2, 12).Value = 2605.9  # L2: was 2838.4443
$ws.Cells.Item(2, 13).Value = -960.7141999999999  # M2: was -949.8620000000001
$ws.Cells.Item(2, 14).Value = -2831.9  # N2: was -3064.4443

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 369462.44  # H32: was 2539.4
$ws.Cells.Item(32, 9).Value = 2550.5693  # I32: was 2009.8948
$ws.Cells.Item(32, 10).Value = 4772404.5  # J32: was 12600
$ws.Cells.Item(32, 11).Value = 2550.5693  # K32: was 2009.8948
$ws.Cells.Item(32, 12).Value = 4772404.5  # L32: was 12600
$ws.Cells.Item(32, 13).Value = -2263.5693  # M32: was -1722.8948
$ws.Cells.Item(32, 14).Value = -4772978.5  # N32: was -13174

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1089.4878  # H61: was 794.6212
$ws.Cells.Item(61, 9).Value = 651.069  # I61: was 499.1887
$ws.Cells.Item(61, 10).Value = 2149  # J61: was 1999.0769
$ws.Cells.Item(61, 11).Value = 651.069  # K61: was 499.1887
$ws.Cells.Item(61, 12).Value = 2149  # L61: was 1999.0769
$ws.Cells.Item(61, 13).Value = -439.069  # M61: was -287.1887
$ws.Cells.Item(61, 14).Value = -2573  # N61: was -2423.0769

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1476.921  # H116: was 1483.3948
$ws.Cells.Item(116, 9).Value = 1073.7142  # I116: was 1062.862
$ws.Cells.Item(116, 10).Value = 2605.9  # J116: was 2838.4443
$ws.Cells.Item(116, 11).Value = 1073.7142  # K116: was 1062.862
$ws.Cells.Item(116, 12).Value = 2605.9  # L116: was 2838.4443
$ws.Cells.Item(116, 13).Value = 1220.2858  # M116: was 1231.138
$ws.Cells.Item(116, 14).Value = -7193.9  # N116: was -7426.4443

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 27029168  # H122: was 1914.0256
$ws.Cells.Item(122, 9).Value = 35716684  # I122: was 2045.5667
$ws.Cells.Item(122, 10).Value = 1337.1111  # J122: was 1475.5555
$ws.Cells.Item(122, 11).Value = 107150052  # K122: was 6136.7001
$ws.Cells.Item(122, 12).Value = 4011.3333  # L122: was 4426.666499999999
$ws.Cells.Item(122, 13).Value = -107147602  # M122: was -3686.7001
$ws.Cells.Item(122, 14).Value = -8911.3333  # N122: was -9326.666499999999

# ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value = 0  # H123: was 49980
$ws.Cells.Item(123, 10).Value = 0  # J123: was 49980
$ws.Cells.Item(123, 12).Value = 0  # L123: was 49980
$ws.Cells.Item(123, 14).ClearContents()  # N123: was -59780

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1951.3334  # H132: was 1584
$ws.Cells.Item(132, 9).Value = 1243.75  # I132: was 838.5714
$ws.Cells.Item(132, 11).Value = 3731.25  # K132: was 2515.7142
$ws.Cells.Item(132, 13).Value = -1201.25  # M132: was 14.28579999999965

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1089.4878  # H136: was 794.6212
$ws.Cells.Item(136, 9).Value = 651.069  # I136: was 499.1887
$ws.Cells.Item(136, 10).Value = 2149  # J136: was 1999.0769
$ws.Cells.Item(136, 11).Value = 1953.207  # K136: was 1497.5661
$ws.Cells.Item(136, 12).Value = 6447  # L136: was 5997.2307
$ws.Cells.Item(136, 13).Value = 596.7930000000001  # M136: was 1052.4339
$ws.Cells.Item(136, 14).Value = -11547  # N136: was -11097.2307

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1476.921  # H3: was 1483.3948
$ws.Cells.Item(3, 9).Value = 1073.7142  # I3: was 1062.862
$ws.Cells.Item(3, 10).Value = 2605.9  # J3: was 2838.4443
$ws.Cells.Item(3, 11).Value = 1073.7142  # K3: was 1062.862
$ws.Cells.Item(3, 12).Value = 2605.9  # L3: was 2838.4443
$ws.Cells.Item(3, 13).Value = -959.7141999999999  # M3: was -948.8620000000001
$ws.Cells.Item(3, 14).Value = -2833.9  # N3: was -3066.4443

# BSM row 34
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(34, 8).Value = 0  # H34: was 4000
$ws.Cells.Item(34, 10).Value = 0  # J34: was 4000
$ws.Cells.Item(34, 12).Value = 0  # L34: was 4000
$ws.Cells.Item(34, 14).ClearContents()  # N34: was -4228

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2108.25  # H99: was 2128.875
$ws.Cells.Item(99, 9).Value = 965.125  # I99: was 1006.375
$ws.Cells.Item(99, 11).Value = 965.125  # K99: was 1006.375
$ws.Cells.Item(99, 13).Value = 532.875  # M99: was 491.625

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4729.6313  # H134: was 4508.5
$ws.Cells.Item(134, 9).Value = 942.3889  # I134: was 908.9474
$ws.Cells.Item(134, 11).Value = 2827.1667  # K134: was 2726.8422
$ws.Cells.Item(134, 13).Value = -292.1667000000002  # M134: was -191.8422

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1091.8  # H99: was 1094.4
$ws.Cells.Item(99, 9).Value = 1023.26086  # I99: was 1007.1111
$ws.Cells.Item(99, 11).Value = 1023.26086  # K99: was 1007.1111
$ws.Cells.Item(99, 13).Value = 474.73914  # M99: was 490.8889

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 1091.8  # H126: was 1094.4
$ws.Cells.Item(126, 9).Value = 1023.26086  # I126: was 1007.1111
$ws.Cells.Item(126, 11).Value = 3069.78258  # K126: was 3021.3333
$ws.Cells.Item(126, 13).Value = -599.7825800000001  # M126: was -551.3332999999998

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1540.826  # H132: was 1739.5946
$ws.Cells.Item(132, 9).Value = 1408.9  # I132: was 1632.7826
$ws.Cells.Item(132, 10).Value = 1788.1875  # J132: was 1915.0714
$ws.Cells.Item(132, 11).Value = 4226.700000000001  # K132: was 4898.3478
$ws.Cells.Item(132, 12).Value = 5364.5625  # L132: was 5745.2142
$ws.Cells.Item(132, 13).Value = -1696.700000000001  # M132: was -2368.3478
$ws.Cells.Item(132, 14).Value = -10424.5625  # N132: was -10805.2142

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1260.8049  # H134: was 1333.5405
$ws.Cells.Item(134, 9).Value = 1302.9429  # I134: was 1379.7188
$ws.Cells.Item(134, 10).Value = 1015  # J134: was 1038
$ws.Cells.Item(134, 11).Value = 3908.8287  # K134: was 4139.1564
$ws.Cells.Item(134, 12).Value = 3045  # L134: was 3114
$ws.Cells.Item(134, 13).Value = -1373.8287  # M134: was -1604.1564
$ws.Cells.Item(134, 14).Value = -8115  # N134: was -8184

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 17232.32  # H46: was 28913.572
$ws.Cells.Item(46, 9).Value = 760  # I46: was 393
$ws.Cells.Item(46, 10).Value = 21350.4  # J46: was 33667
$ws.Cells.Item(46, 11).Value = 2280  # K46: was 1179
$ws.Cells.Item(46, 12).Value = 64051.2  # L46: was 101001
$ws.Cells.Item(46, 13).Value = -2189  # M46: was -1088
$ws.Cells.Item(46, 14).Value = -64233.2  # N46: was -101183

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 1842.8572  # H70: was 1354.4445
$ws.Cells.Item(70, 9).Value = 1200  # I70: was 1031.6666
$ws.Cells.Item(70, 10).Value = 1950  # J70: was 2000
$ws.Cells.Item(70, 11).Value = 3600  # K70: was 3094.9998
$ws.Cells.Item(70, 12).Value = 5850  # L70: was 6000
$ws.Cells.Item(70, 13).Value = -3285  # M70: was -2779.9998
$ws.Cells.Item(70, 14).Value = -6480  # N70: was -6630

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(73, 8).Value = 1842.8572  # H73: was 1354.4445
$ws.Cells.Item(73, 9).Value = 1200  # I73: was 1031.6666
$ws.Cells.Item(73, 10).Value = 1950  # J73: was 2000
$ws.Cells.Item(73, 11).Value = 3600  # K73: was 3094.9998
$ws.Cells.Item(73, 12).Value = 5850  # L73: was 6000
$ws.Cells.Item(73, 13).Value = -2508  # M73: was -2002.9998
$ws.Cells.Item(73, 14).Value = -8034  # N73: was -8184

# CUL row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(116, 8).Value = 1291.125  # H116: was 900
$ws.Cells.Item(116, 9).Value = 857.25  # I116: was 900
$ws.Cells.Item(116, 10).Value = 1725  # J116: was 0
$ws.Cells.Item(116, 11).Value = 2571.75  # K116: was 2700
$ws.Cells.Item(116, 12).Value = 5175  # L116: was 0
$ws.Cells.Item(116, 13).Value = 870.25  # M116: was 742
$ws.Cells.Item(116, 14).Value = -12059  # N116: was None

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 632.3913  # H122: was 626.3043
$ws.Cells.Item(122, 9).Value = 414.2857  # I122: was 372.22223
$ws.Cells.Item(122, 10).Value = 727.8125  # J122: was 789.6429000000001
$ws.Cells.Item(122, 11).Value = 3728.5713  # K122: was 3350.00007
$ws.Cells.Item(122, 12).Value = 6550.3125  # L122: was 7106.7861
$ws.Cells.Item(122, 13).Value = -1278.5713  # M122: was -900.0000700000001
$ws.Cells.Item(122, 14).Value = -11450.3125  # N122: was -12006.7861

# GSM row 20
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 0  # H20: was 1999.2222
$ws.Cells.Item(20, 10).Value = 0  # J20: was 1999.2222
$ws.Cells.Item(20, 12).Value = 0  # L20: was 1999.2222
$ws.Cells.Item(20, 14).ClearContents()  # N20: was -2489.2222

# GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 2000  # H24: was 4007
$ws.Cells.Item(24, 10).Value = 2000  # J24: was 4007
$ws.Cells.Item(24, 12).Value = 2000  # L24: was 4007
$ws.Cells.Item(24, 14).Value = -2346  # N24: was -4353

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2013.7084  # H122: was 1824.3793
$ws.Cells.Item(122, 9).Value = 1851.3125  # I122: was 1743.9445
$ws.Cells.Item(122, 10).Value = 2338.5  # J122: was 1956
$ws.Cells.Item(122, 11).Value = 5553.9375  # K122: was 5231.833500000001
$ws.Cells.Item(122, 12).Value = 7015.5  # L122: was 5868
$ws.Cells.Item(122, 13).Value = -3103.9375  # M122: was -2781.833500000001
$ws.Cells.Item(122, 14).Value = -11915.5  # N122: was -10768

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1428.909  # H132: was 1313.2545
$ws.Cells.Item(132, 9).Value = 1388.8572  # I132: was 1294.159
$ws.Cells.Item(132, 10).Value = 1584.6666  # J132: was 1389.6364
$ws.Cells.Item(132, 11).Value = 4166.571599999999  # K132: was 3882.477
$ws.Cells.Item(132, 12).Value = 4753.9998  # L132: was 4168.9092
$ws.Cells.Item(132, 13).Value = -1636.571599999999  # M132: was -1352.477
$ws.Cells.Item(132, 14).Value = -9813.9998  # N132: was -9228.9092

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2038.1818  # H40: was 2020.129
$ws.Cells.Item(40, 9).Value = 1796.6666  # I40: was 1804.9333
$ws.Cells.Item(40, 10).Value = 2328  # J40: was 2221.875
$ws.Cells.Item(40, 11).Value = 1796.6666  # K40: was 1804.9333
$ws.Cells.Item(40, 12).Value = 2328  # L40: was 2221.875
$ws.Cells.Item(40, 13).Value = -1660.6666  # M40: was -1668.9333
$ws.Cells.Item(40, 14).Value = -2600  # N40: was -2493.875

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48, 8).Value = 5291.5713  # H48: was 5800
$ws.Cells.Item(48, 9).Value = 4020.5  # I48: was 0
$ws.Cells.Item(48, 11).Value = 4020.5  # K48: was 0
$ws.Cells.Item(48, 13).Value = -3359.5  # M48: was None

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2179.087  # H122: was 2039.3462
$ws.Cells.Item(122, 9).Value = 2048.2354  # I122: was 1886.2
$ws.Cells.Item(122, 11).Value = 6144.706200000001  # K122: was 5658.6
$ws.Cells.Item(122, 13).Value = -3694.706200000001  # M122: was -3208.6

# WVR row 45
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 12706.375  # H45: was 12253.125
$ws.Cells.Item(45, 9).Value = 0  # I45: was 5000
$ws.Cells.Item(45, 10).Value = 12706.375  # J45: was 13289.286
$ws.Cells.Item(45, 11).Value = 0  # K45: was 5000
$ws.Cells.Item(45, 12).Value = 12706.375  # L45: was 13289.286
$ws.Cells.Item(45, 13).ClearContents()  # M45: was -4509
$ws.Cells.Item(45, 14).Value = -13688.375  # N45: was -14271.286

# WVR row 74
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 36182.332  # H74: was 12820.833
$ws.Cells.Item(74, 9).Value = 5569  # I74: was 0
$ws.Cells.Item(74, 10).Value = 42305  # J74: was 12820.833
$ws.Cells.Item(74, 11).Value = 5569  # K74: was 0
$ws.Cells.Item(74, 12).Value = 42305  # L74: was 12820.833
$ws.Cells.Item(74, 13).Value = -4633  # M74: was None
$ws.Cells.Item(74, 14).Value = -44177  # N74: was -14692.833

# WVR row 77
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(77, 8).Value = 36182.332  # H77: was 12820.833
$ws.Cells.Item(77, 9).Value = 5569  # I77: was 0
$ws.Cells.Item(77, 10).Value = 42305  # J77: was 12820.833
$ws.Cells.Item(77, 11).Value = 16707  # K77: was 0
$ws.Cells.Item(77, 12).Value = 126915  # L77: was 38462.499
$ws.Cells.Item(77, 13).Value = -12027  # M77: was None
$ws.Cells.Item(77, 14).Value = -136275  # N77: was -47822.499

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1469.2307  # H122: was 1413.8966
$ws.Cells.Item(122, 9).Value = 1395.0588  # I122: was 1332.421
$ws.Cells.Item(122, 10).Value = 1609.3334  # J122: was 1568.7
$ws.Cells.Item(122, 11).Value = 4185.1764  # K122: was 3997.263
$ws.Cells.Item(122, 12).Value = 4828.0002  # L122: was 4706.1
$ws.Cells.Item(122, 13).Value = -1735.1764  # M122: was -1547.263
$ws.Cells.Item(122, 14).Value = -9728.0002  # N122: was -9606.1
